$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "backpack" synthesis column (H): header + recipe, formatted the same
# way as the adjacent "黃金鏟子" column (G) so it reuses that cell style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "黃金十字鎬"

$ws.Range("H2").Value = " 樹枝*2,燧石*2,黃金*1"

# Give the new column a width matching its neighbours.
$ws.Columns.Item(8).ColumnWidth = 26.72712053571429

# Leave the view scrolled/selected on the new column, like the author did.
[void]$ws.Range("I1:I1048576").Select()
